$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 2 for the latest 4D Box draw.
# This shifts the previous rows 2-39 down to rows 3-40.
$ws.Rows.Item(2).Insert()

# Insert() copies the header row's (yellow/bold) formatting down onto the new
# row; strip that back off *before* writing any values (clearing formats
# after writing the checkmark/emoji text can make the engine pick a
# substitute font for it, which we don't want).
$ws.Range("A2:C2").ClearFormats()

# --- Row 2: new latest draw entry (date unchanged, new numbers/result) ---
$ws.Range("A2").Value = "2/7/2025 (Wed)"
$ws.Range("B2").Value = "6 4 1 8" + [char]10 + "2 5 2 6" + [char]10 + "7 9 9 7" + [char]10 + "5 6 3 0"
$ws.Range("C2").Value = "✅ Direct: 14/4042 (0.35%)" + [char]10 + "✅ iBet: 14/216 (6.48%)"

# Re-apply the plain wrap-text style used by all the other data rows in the
# 4D Box / Result columns.
$ws.Range("B2:C2").WrapText = $true

# Remove the auto row-height bump caused by the multi-line text so row 2
# keeps the sheet's default (non custom) height, matching the source pattern.
$ws.Rows.Item(2).EntireRow.AutoFit()

# --- Row 3 now holds what used to be row 2's data; give it the same
# 60pt custom row height used by all the other historical data rows. ---
$ws.Rows.Item(3).RowHeight = 60

# The insert pushed one extra blank placeholder row (which used to have the
# 60pt custom height) into what is now row 14; restore it to the sheet's
# default (auto) height to match the tail formatting pattern.
$ws.Rows.Item(14).EntireRow.AutoFit()

Write-Host "done"
